$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name) 부서별원가분석 -> 부서별원가
$ws.Name = "부서별원가"

# Insert a new column before column A (shifts existing data right: old A->B, ... old F->G)
$ws.Columns.Item(1).Insert()

# Remove the now-shifted last column (was G:작업시간, now column H) entirely - not present in target
$ws.Columns.Item(8).Delete()

# Remove rows 5 through 8 (only 3 data rows + header remain in target)
$ws.Range("A5:G8").EntireRow.Delete()

# --- Header row ---
$ws.Range("A1").Value = "년월"
$ws.Range("B1").Value = "부서코드"
$ws.Range("C1").Value = "부서명"
$ws.Range("D1").Value = "제품코드"
$ws.Range("E1").Value = "제품명"
$ws.Range("F1").Value = "원가"
$ws.Range("G1").Value = "수량"

# --- Row 2 ---
$ws.Range("A2").Value = "202411"
$ws.Range("B2").Value = "D001"
$ws.Range("C2").Value = "생산1팀"
$ws.Range("D2").Value = "MDL-001"
$ws.Range("E2").Value = "디스플레이 패널"
$ws.Range("F2").Value = 15000
$ws.Range("G2").Value = 100

# --- Row 3 ---
$ws.Range("A3").Value = "202411"
$ws.Range("B3").Value = "D002"
$ws.Range("C3").Value = "생산2팀"
$ws.Range("D3").Value = "MDL-002"
$ws.Range("E3").Value = "LED 모듈"
$ws.Range("F3").Value = 8000
$ws.Range("G3").Value = 200

# --- Row 4 ---
$ws.Range("A4").Value = "202411"
$ws.Range("B4").Value = "D003"
$ws.Range("C4").Value = "품질관리팀"
$ws.Range("D4").Value = "MDL-003"
$ws.Range("E4").Value = "컨트롤러"
$ws.Range("F4").Value = 12000
$ws.Range("G4").Value = 150

Write-Host "done"
